# Stop tracking INVESTMENT excel file
# Refresh the daily snapshot values in "CURRENT STATUS" sheet (rows 2-17)
# Date (col B), Close (col C), C_EMA1 (col F), C_EMA2 (col G) and
# C_RSI_V (col J) are recalculated for the new as-of date, and a few
# C_Pos (col L) signals flip to "B_N_Y".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - ICIGOL
$ws.Range("B2").Value = 45895
$ws.Range("C2").Value = 86.59
$ws.Range("F2").Value = 85.48114408005708
$ws.Range("G2").Value = 84.75881969628711
$ws.Range("J2").Value = 59.17351827644686

# Row 3 - ICIPSE
$ws.Range("B3").Value = 45895
$ws.Range("C3").Value = 116.33
$ws.Range("F3").Value = 114.8822831169633
$ws.Range("G3").Value = 101.4074082968453
$ws.Range("J3").Value = 65.54654453288495

# Row 4 - KOTN50
$ws.Range("B4").Value = 45895
$ws.Range("C4").Value = 48.25
$ws.Range("F4").Value = 48.92811601496058
$ws.Range("G4").Value = 49.37356094495915
$ws.Range("J4").Value = 42.90347922663884

# Row 5 - SBINIF
$ws.Range("B5").Value = 45895
$ws.Range("C5").Value = 263.99
$ws.Range("F5").Value = 264.9394782148985
$ws.Range("G5").Value = 264.7290311795473
$ws.Range("J5").Value = 47.30924099333645

# Row 6 - CPSETF
$ws.Range("B6").Value = 45895
$ws.Range("C6").Value = 88.56
$ws.Range("F6").Value = 89.63748989688654
$ws.Range("G6").Value = 90.87973945871273
$ws.Range("J6").Value = 32.06904172280153

# Row 7 - HDF250
$ws.Range("B7").Value = 45895
$ws.Range("C7").Value = 168.77
$ws.Range("F7").Value = 170.7685201080272
$ws.Range("G7").Value = 171.9433314873023
$ws.Range("J7").Value = 46.37145806981436

# Row 8 - ICI150 (C_Pos flips S -> B_N_Y)
$ws.Range("B8").Value = 45895
$ws.Range("C8").Value = 21.79
$ws.Range("F8").Value = 21.97747013795918
$ws.Range("G8").Value = 21.9092542656204
$ws.Range("J8").Value = 48.23490240805658
$ws.Range("L8").Value = "B_N_Y"

# Row 9 - ICIAUT
$ws.Range("B9").Value = 45895
$ws.Range("C9").Value = 26.01
$ws.Range("F9").Value = 25.94873099212012
$ws.Range("G9").Value = 25.65535626304575
$ws.Range("J9").Value = 70.39893240239451

# Row 10 - ICIFMC
$ws.Range("B10").Value = 45895
$ws.Range("C10").Value = 60.04
$ws.Range("F10").Value = 59.50702631040837
$ws.Range("G10").Value = 59.07986340262906
$ws.Range("J10").Value = 54.49285588394612

# Row 11 - ICIHEA (C_Pos flips B -> B_N_Y)
$ws.Range("B11").Value = 45895
$ws.Range("C11").Value = 148.69
$ws.Range("F11").Value = 149.5348902227552
$ws.Range("G11").Value = 146.0244058914295
$ws.Range("J11").Value = 48.35756822698924
$ws.Range("L11").Value = "B_N_Y"

# Row 12 - MOTNAS
$ws.Range("B12").Value = 45895
$ws.Range("C12").Value = 202.49
$ws.Range("F12").Value = 202.8187961612267
$ws.Range("G12").Value = 187.5879599890849
$ws.Range("J12").Value = 58.34945619049549

# Row 13 - SBIEIT
$ws.Range("B13").Value = 45895
$ws.Range("C13").Value = 393
$ws.Range("F13").Value = 389.9314420348093
$ws.Range("G13").Value = 399.5858584837641
$ws.Range("J13").Value = 63.33395588756564

# Row 14 - ICIEXB (C_Pos flips S -> B_N_Y)
$ws.Range("B14").Value = 45895
$ws.Range("C14").Value = 29.38
$ws.Range("F14").Value = 29.91288827469454
$ws.Range("G14").Value = 29.87905319256971
$ws.Range("J14").Value = 44.57734174937318
$ws.Range("L14").Value = "B_N_Y"

# Row 15 - SBIBAN
$ws.Range("B15").Value = 45895
$ws.Range("C15").Value = 558.84
$ws.Range("F15").Value = 566.9628495599803
$ws.Range("G15").Value = 570.2580523183403
$ws.Range("J15").Value = 40.28750617781254

# Row 16 - MAHMAH
$ws.Range("B16").Value = 45895
$ws.Range("C16").Value = 3330.8
$ws.Range("F16").Value = 3180.867321256303
$ws.Range("G16").Value = 3012.119749076148
$ws.Range("J16").Value = 51.76562365762512

# Row 17 - GODPRO
$ws.Range("B17").Value = 45895
$ws.Range("C17").Value = 2007.1
$ws.Range("F17").Value = 2072.471307808276
$ws.Range("G17").Value = 2211.316987214698
$ws.Range("J17").Value = 40.03153029554537
